$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null

$ws.Range("H97").Value = 6363.25
$ws.Range("J97").Value = 6363.25
$ws.Range("L97").Value = 19089.75
$ws.Range("N97").Value = -20081.75

$ws.Range("H112").Value = 8775.412
$ws.Range("I112").Value = 966.3333
$ws.Range("J112").Value = 10448.786
$ws.Range("K112").Value = 2898.9999
$ws.Range("L112").Value = 31346.358
$ws.Range("M112").Value = -1790.9999
$ws.Range("N112").Value = -33562.358

$ws.Range("H113").Value = 8129.3335
$ws.Range("I113").Value = 7750
$ws.Range("K113").Value = 7750
$ws.Range("M113").Value = -4496

$ws.Range("H138").Value = 1696763.9
$ws.Range("J138").Value = 3000173
$ws.Range("L138").Value = 9000519
$ws.Range("N138").Value = -9010799

$ws.Range("H141").Value = 2203.4167
$ws.Range("I141").Value = 2203.4167
$ws.Range("K141").Value = 6610.250100000001
$ws.Range("M141").Value = -1430.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18103.297
$ws.Range("I32").Value = 18565.213
$ws.Range("K32").Value = 18565.213
$ws.Range("M32").Value = -18278.213

$ws.Range("H74").Value = 3243.0364
$ws.Range("I74").Value = 1188.186
$ws.Range("K74").Value = 1188.186
$ws.Range("M74").Value = -314.1859999999999

$ws.Range("H77").Value = 3243.0364
$ws.Range("I77").Value = 1188.186
$ws.Range("K77").Value = 5940.929999999999
$ws.Range("M77").Value = -1572.929999999999

$ws.Range("H94").Value = 75000
$ws.Range("J94").Value = 75000
$ws.Range("L94").Value = 75000
$ws.Range("N94").Value = -76802

$ws.Range("H97").Value = 712
$ws.Range("I97").Value = 522.36365
$ws.Range("J97").Value = 1407.3334
$ws.Range("K97").Value = 522.36365
$ws.Range("L97").Value = 1407.3334
$ws.Range("M97").Value = -26.36365000000001
$ws.Range("N97").Value = -2399.3334

$ws.Range("H112").Value = 31571.143
$ws.Range("J112").Value = 31571.143
$ws.Range("L112").Value = 31571.143
$ws.Range("N112").Value = -34525.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 584.7222
$ws.Range("J80").Value = 315.625
$ws.Range("L80").Value = 315.625
$ws.Range("N80").Value = -2311.625

$ws.Range("H83").Value = 584.7222
$ws.Range("J83").Value = 315.625
$ws.Range("L83").Value = 1578.125
$ws.Range("N83").Value = -11562.125

$ws.Range("H107").Value = 776.4761999999999
$ws.Range("I107").Value = 630.3333
$ws.Range("J107").Value = 1653.3334
$ws.Range("K107").Value = 630.3333
$ws.Range("L107").Value = 1653.3334
$ws.Range("M107").Value = 1289.6667
$ws.Range("N107").Value = -5493.3334

$ws.Range("H133").Value = 68500
$ws.Range("I133").Value = 68500
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 68500
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -63440
$ws.Range("N133").Value = $null

$ws.Range("H134").Value = 3466.1177
$ws.Range("I134").Value = 3398.3333
$ws.Range("K134").Value = 10194.9999
$ws.Range("M134").Value = -7659.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1963350.4
$ws.Range("I31").Value = 3032410.2
$ws.Range("K31").Value = 3032410.2
$ws.Range("M31").Value = -3032115.2

$ws.Range("H34").Value = 1963350.4
$ws.Range("I34").Value = 3032410.2
$ws.Range("K34").Value = 3032410.2
$ws.Range("M34").Value = -3032208.2

$ws.Range("H62").Value = 5108.1665
$ws.Range("I62").Value = 4949.5
$ws.Range("K62").Value = 4949.5
$ws.Range("M62").Value = -4325.5

$ws.Range("H65").Value = 5108.1665
$ws.Range("I65").Value = 4949.5
$ws.Range("K65").Value = 24747.5
$ws.Range("M65").Value = -21627.5

$ws.Range("H98").Value = 88700
$ws.Range("J98").Value = 88700
$ws.Range("L98").Value = 88700
$ws.Range("N98").Value = -93192

$ws.Range("H107").Value = 415.66666
$ws.Range("I107").Value = 253.45
$ws.Range("K107").Value = 253.45
$ws.Range("M107").Value = 1666.55

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = $null

$ws.Range("H132").Value = 2447.2307
$ws.Range("I132").Value = 1983.6364
$ws.Range("K132").Value = 5950.9092
$ws.Range("M132").Value = -3420.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 443.75
$ws.Range("I2").Value = 475
$ws.Range("K2").Value = 2850
$ws.Range("M2").Value = -2737

$ws.Range("H68").Value = 4372.609
$ws.Range("I68").Value = 1926.3334
$ws.Range("J68").Value = 4739.55
$ws.Range("K68").Value = 5779.0002
$ws.Range("L68").Value = 14218.65
$ws.Range("M68").Value = -4968.0002
$ws.Range("N68").Value = -15840.65

$ws.Range("H71").Value = 4372.609
$ws.Range("I71").Value = 1926.3334
$ws.Range("J71").Value = 4739.55
$ws.Range("K71").Value = 17337.0006
$ws.Range("L71").Value = 42655.95
$ws.Range("M71").Value = -13281.0006
$ws.Range("N71").Value = -50767.95

$ws.Range("H113").Value = 758.6667
$ws.Range("J113").Value = 836.1
$ws.Range("L113").Value = 2508.3
$ws.Range("N113").Value = -6848.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 436.63333
$ws.Range("I2").Value = 667.7143
$ws.Range("K2").Value = 667.7143
$ws.Range("M2").Value = -554.7143

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

$ws.Range("H97").Value = 824.0476
$ws.Range("I97").Value = 647.35297
$ws.Range("K97").Value = 647.35297
$ws.Range("M97").Value = -151.35297

$ws.Range("H122").Value = 166670000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 166670000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 500010000
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -500014900

$ws.Range("H132").Value = 10156.23
$ws.Range("I132").Value = 9636.637000000001
$ws.Range("J132").Value = 13014
$ws.Range("K132").Value = 28909.911
$ws.Range("L132").Value = 39042
$ws.Range("M132").Value = -26379.911
$ws.Range("N132").Value = -44102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2750.7693
$ws.Range("I7").Value = 2497.611
$ws.Range("K7").Value = 2497.611
$ws.Range("M7").Value = -2385.611

$ws.Range("H68").Value = 1999.2
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 1999.5
$ws.Range("K68").Value = 1999
$ws.Range("L68").Value = 1999.5
$ws.Range("M68").Value = -1250
$ws.Range("N68").Value = -3497.5

$ws.Range("H71").Value = 1999.2
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 1999.5
$ws.Range("K71").Value = 9995
$ws.Range("L71").Value = 9997.5
$ws.Range("M71").Value = -6251
$ws.Range("N71").Value = -17485.5

$ws.Range("H93").Value = 1262
$ws.Range("J93").Value = 1194.6666
$ws.Range("L93").Value = 1194.6666
$ws.Range("N93").Value = -3690.6666

$ws.Range("H126").Value = 2750.7693
$ws.Range("I126").Value = 2497.611
$ws.Range("K126").Value = 7492.833
$ws.Range("M126").Value = -5022.833

$ws.Range("H132").Value = 4519.4517
$ws.Range("I132").Value = 3911.7693
$ws.Range("J132").Value = 4958.3335
$ws.Range("K132").Value = 11735.3079
$ws.Range("L132").Value = 14875.0005
$ws.Range("M132").Value = -9205.3079
$ws.Range("N132").Value = -19935.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13230.125
$ws.Range("J41").Value = 13120.143
$ws.Range("L41").Value = 13120.143
$ws.Range("N41").Value = -13900.143

$ws.Range("H62").Value = 19498.25
$ws.Range("J62").Value = 19498.25
$ws.Range("L62").Value = 19498.25
$ws.Range("N62").Value = -20746.25

$ws.Range("H65").Value = 19498.25
$ws.Range("J65").Value = 19498.25
$ws.Range("L65").Value = 97491.25
$ws.Range("N65").Value = -103731.25

$ws.Range("H81").Value = 4458.7295
$ws.Range("I81").Value = 4169.3438
$ws.Range("K81").Value = 8338.687599999999
$ws.Range("M81").Value = -7277.687599999999

$ws.Range("H84").Value = 4458.7295
$ws.Range("I84").Value = 4169.3438
$ws.Range("K84").Value = 41693.43799999999
$ws.Range("M84").Value = -36389.43799999999

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

$ws.Range("H132").Value = 18787
$ws.Range("I132").Value = 21113.084
$ws.Range("K132").Value = 63339.25199999999
$ws.Range("M132").Value = -60809.25199999999

$ws.Range("H140").Value = 64966.668
$ws.Range("I140").Value = 64900
$ws.Range("K140").Value = 64900
$ws.Range("M140").Value = -59720
